# Update result values on each year sheet (row 2) with the latest server results.
$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 973.9537847600009
$ws2025.Range("E2").Value = 28982.37596598056
$ws2025.Range("I2").Value = 16175.28135478
$ws2025.Range("L2").Value = 48524.529503538
$ws2025.Range("M2").Value = 10590.587968015
$ws2025.Range("N2").Value = 7155.07579047334
$ws2025.Range("O2").Value = 6980.325566461758

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 5712.560177842886
$ws2030.Range("E2").Value = 56106.05588781912
$ws2030.Range("I2").Value = 44217.8984721661
$ws2030.Range("L2").Value = 66966.57749858923
$ws2030.Range("M2").Value = 21984.28023276101
$ws2030.Range("N2").Value = 10593.94009226292
$ws2030.Range("O2").Value = 12063.05946690077

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 2861.961401238371
$ws2035.Range("B2").Value = 8026.889663087295
$ws2035.Range("E2").Value = 67297.73995507321
$ws2035.Range("I2").Value = 59256.42575923612
$ws2035.Range("L2").Value = 66966.57749858923
$ws2035.Range("M2").Value = 25464.6214365565
$ws2035.Range("N2").Value = 15131.8804243991
$ws2035.Range("O2").Value = 14761.05157597961

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Value = 2861.961401238371
$ws2040.Range("B2").Value = 8026.889663087295
$ws2040.Range("E2").Value = 67297.73995507321
$ws2040.Range("I2").Value = 59256.42575923612
$ws2040.Range("L2").Value = 66966.57749858923
$ws2040.Range("M2").Value = 25464.6214365565
$ws2040.Range("N2").Value = 15236.2891206102
$ws2040.Range("O2").Value = 14761.05157597961

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 6302.873118834019
$ws2045.Range("B2").Value = 8026.889663087295
$ws2045.Range("E2").Value = 67297.73995507321
$ws2045.Range("I2").Value = 59256.42575923612
$ws2045.Range("L2").Value = 66966.57749858923
$ws2045.Range("M2").Value = 25464.6214365565
$ws2045.Range("N2").Value = 15771.80030518411
$ws2045.Range("O2").Value = 17096.51756232827

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Value = 6302.873118834019
$ws2050.Range("B2").Value = 8026.889663087295
$ws2050.Range("E2").Value = 67297.73995507321
$ws2050.Range("I2").Value = 59256.42575923612
$ws2050.Range("L2").Value = 66966.57749858923
$ws2050.Range("M2").Value = 25464.6214365565
$ws2050.Range("N2").Value = 15771.80030518411
$ws2050.Range("O2").Value = 17096.51756232827
